# Develop: reset playerModel structure and calculate the total number of each attribute
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..V correspond to attribute stats for each player row.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

# Capture row 3 (André Onana) values before resetting it, then move them to
# row 4 (Altay Bayindir), and zero out row 3 ("reset playerModel structure").
$row3Values = @{}
foreach ($col in $cols) {
    $row3Values[$col] = $ws.Range("$col`3").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col`4").Value2 = $row3Values[$col]
    $ws.Range("$col`3").Value2 = 0
}

# Update the active selection / view to match the new state.
$ws.Range("P9").Select()
